$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 15 (shifts old rows 15-22 down to 16-23)
$ws.Rows.Item(15).Insert()

# New to-do item under "Destinations": "Optional: drop down menu"
$ws.Range("H14").Value = "Optional: drop down menu"

# New content at the bottom of the list: "Overall" / "Sort CSS"
$ws.Range("F25").Value = "Overall"
$ws.Range("H25").Value = "Sort CSS"

# Update the view: scroll/selection to reflect the new active cell
$ws.Application.GoTo($ws.Range("A3"))
$ws.Range("H14").Select()

# Add a second worksheet (after Sheet1) with important notes
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "important notes"
$ws2.Range("B3").Value = "12 columns = 1200 px"
$ws2.Range("B5").Select()

$ws.Select()
